$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Heading2" style from the section-heading paragraphs so they
#    become normal paragraphs (matches removal of <w:pPr><w:pStyle .../></w:pPr>).
# ---------------------------------------------------------------------------
$headingTexts = @(
    "Introduction",
    "The Transcendent Nature of God",
    "Historical Revelations of God",
    "Moral Foundations and Expectations",
    "Rituals, Practices, and Legal Structures",
    "Conclusion"
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd("`r", "`n", [char]7)
    if ($headingTexts -contains $text) {
        $p.Style = "Normal"
    }
}

# ---------------------------------------------------------------------------
# 2) Replace the inline citation markers with the new reference tags.
#    Each paragraph gets its citations collapsed onto a single reference id,
#    so the replacement is scoped per-paragraph (same author name appears in
#    different paragraphs but must map to different ids there).
# ---------------------------------------------------------------------------
function Replace-InParagraphRange {
    param($range, [string]$find, [string]$replace)
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text

    if ($text.StartsWith("In Judaism, the transcendent nature of God")) {
        Replace-InParagraphRange $p.Range "(Green)" "(Ref-f239155)"
        Replace-InParagraphRange $p.Range "(Volli)" "(Ref-f239155)"
    }
    elseif ($text.StartsWith("Additionally, the immanent role of God")) {
        Replace-InParagraphRange $p.Range "(Eisenberg)" "(Ref-f706226)"
        Replace-InParagraphRange $p.Range "(Volli)" "(Ref-f706226)"
    }
    elseif ($text.StartsWith("Throughout history, God's revelations")) {
        Replace-InParagraphRange $p.Range "(Eisenberg)" "(Johnson)"
        Replace-InParagraphRange $p.Range "(Volli)" "(Johnson)"
    }
    elseif ($text.StartsWith("The moral principles and expectations")) {
        Replace-InParagraphRange $p.Range "(Eisenberg)" "(Ref-s005380)"
        Replace-InParagraphRange $p.Range "(Green)" "(Ref-s005380)"
    }
    elseif ($text.StartsWith("Jewish rituals and practices play a pivotal role")) {
        Replace-InParagraphRange $p.Range "(Sosis)" "(Ref-f097121)"
        Replace-InParagraphRange $p.Range "(Kelley et al.)" "(Ref-f097121)"
    }
}

Write-Host "edit complete"
